$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

# Add new row 10 entry: APA Site Visit at Utah State University / School Psychology PhD Program Accreditation, Nov 2025
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "APA Site Visit "
$ws.Cells.Item(10, 3).Value = "School Psychology PhD Program Accreditation"
$ws.Cells.Item(10, 4).Value = "Nov 2025"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "Utah State University"

$ws.Range("D10").Select()
